$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# Fill data columns: I is constant 1, J duplicates column H, for rows 2-37.
for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $hVal = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 10).Value = $hVal
}
